$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (bump by one day)
$ws.Range("A1").Value = 45311

# Update the price list values in column D
$ws.Range("D33").Value = 351
$ws.Range("D34").Value = 431.1
$ws.Range("D35").Value = 460
$ws.Range("D36").Value = 580
$ws.Range("D37").Value = 731
$ws.Range("D38").Value = 584.3
$ws.Range("D39").Value = 715
$ws.Range("D40").Value = 890.3
$ws.Range("D41").Value = 1034
